$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-08 Monday", "2024-07-09 Tuesday"),
    @("908÷7=129, 5", "613÷7=87, 4"),
    @("886÷4=221, 2", "720÷9=80, 0"),
    @("127÷7=18, 1", "730÷2=365, 0"),
    @("616÷8=77, 0", "383÷7=54, 5"),
    @("558÷4=139, 2", "331÷3=110, 1"),
    @("149÷2=74, 1", "789÷8=98, 5"),
    @("583÷7=83, 2", "279÷4=69, 3"),
    @("994÷3=331, 1", "383÷2=191, 1"),
    @("645÷8=80, 5", "369÷3=123, 0"),
    @("489÷3=163, 0", "760÷8=95, 0"),
    @("671÷5=134, 1", "231÷2=115, 1"),
    @("281÷6=46, 5", "397÷6=66, 1"),
    @("687÷7=98, 1", "549÷7=78, 3"),
    @("913÷9=101, 4", "857÷4=214, 1"),
    @("899÷6=149, 5", "808÷6=134, 4"),
    @("614÷6=102, 2", "512÷7=73, 1"),
    @("129÷2=64, 1", "107÷8=13, 3"),
    @("334÷3=111, 1", "490÷6=81, 4"),
    @("384÷7=54, 6", "320÷8=40, 0"),
    @("164÷9=18, 2", "372÷6=62, 0"),
    @("494÷5=98, 4", "976÷5=195, 1"),
    @("598÷7=85, 3", "708÷4=177, 0"),
    @("416÷4=104, 0", "494÷3=164, 2"),
    @("137÷5=27, 2", "722÷7=103, 1"),
    @("628÷7=89, 5", "442÷9=49, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}
